# MAI_holdings.xlsx update script
# - Updates the "as of" date in the confidential disclaimer text (A10)
# - Updates weight/percent-change figures in D2:E7
#
# The worksheet is protected, so we must unprotect it before writing values
# and re-protect it afterwards to restore the original protected state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wasProtected = $ws.ProtectContents

if ($wasProtected) {
    $ws.Unprotect()
}

# Update the confidential disclaimer date from 2021-06-09 to 2021-06-10
$oldText = $ws.Range("A10").Value2
$newText = $oldText -replace "2021-06-09", "2021-06-10"
$ws.Range("A10").Value = $newText

# Update the weight (D) and percent change (E) values for rows 2-6
$ws.Range("D2").Value = 0.4912443332613747
$ws.Range("E2").Value = 0.0007754943776658862

$ws.Range("D3").Value = 0.330228229655524
$ws.Range("E3").Value = 0.001680986178558053

$ws.Range("D4").Value = 0.09108203339241193
$ws.Range("E4").Value = 0.002821576763485378

$ws.Range("D5").Value = 0.05815523583368254
$ws.Range("E5").Value = 0.00183003545693694

$ws.Range("D6").Value = 0.0292901678570068
$ws.Range("E6").Value = 0.02603327965646796

# Row 7 ("Total") only has the percent-change value updated
$ws.Range("E7").Value = 0.002062006531904581

if ($wasProtected) {
    $ws.Protect()
}
